$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.803.01"
$ws.Range("E2").Value = "  -3.92%  "
$ws.Range("D3").Value = "3.313.37"
$ws.Range("E3").Value = "  -5.64%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "183.23"
$ws.Range("E5").Value = "  -7.56%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "531.48"
$ws.Range("E6").Value = "  -3.34%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.605"
$ws.Range("E7").Value = "  +0.48%  "
$ws.Range("D8").Value = "3.306.04"
$ws.Range("E8").Value = "  -5.68%  "
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("E10").Value = "  -5.06%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "60.34"
$ws.Range("E11").Value = "  -4.20%  "
$ws.Range("E12").Value = "  -5.94%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000262"
$ws.Range("E13").Value = "  -2.31%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.15"
$ws.Range("E14").Value = "  -6.63%  "
$ws.Range("D15").Value = "3.849.47"
$ws.Range("E15").Value = "  -5.00%  "
$ws.Range("D16").Value = "3.318.21"
$ws.Range("E16").Value = "  -5.00%  "
$ws.Range("E17").Value = "  -4.64%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.78"
$ws.Range("E18").Value = "  -3.45%  "
$ws.Range("D19").Value = "64.677.88"
$ws.Range("E19").Value = "  -3.61%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.22"
$ws.Range("E20").Value = "  -4.96%  "
$ws.Range("E21").Value = "  -6.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "377.32"
$ws.Range("E22").Value = "  -3.30%  "
$ws.Range("E23").Value = "  -4.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.51"
$ws.Range("E24").Value = "  -0.86%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.19"
$ws.Range("E25").Value = "  -5.49%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.89"
$ws.Range("E26").Value = "  +3.45%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.09"
$ws.Range("E27").Value = "  -1.19%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.70"
$ws.Range("E28").Value = "  -3.15%  "
$ws.Range("E29").Value = "  -4.29%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.49"
$ws.Range("E30").Value = "  -2.79%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "29.10"
$ws.Range("E31").Value = "  -5.99%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "653.01"
$ws.Range("E32").Value = "  -3.29%  "
$ws.Range("E33").Value = "  -2.95%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.37"
$ws.Range("E34").Value = "  -2.66%  "
$ws.Range("E35").Value = "  -3.42%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "59.42"
$ws.Range("E36").Value = "  -5.98%  "
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.394"
$ws.Range("E38").Value = "  -1.20%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.71"
$ws.Range("E39").Value = "  -5.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  +0.30%  "
$ws.Range("D41").Value = "0.0₃0714"
$ws.Range("E41").Value = "  +5.81%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.127"
$ws.Range("E42").Value = "  -2.36%  "
$ws.Range("D43").Value = "2.898.09"
$ws.Range("E43").Value = "  -5.41%  "
$ws.Range("E44").Value = "  +0.41%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.71"
$ws.Range("E45").Value = "  -9.31%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0403"
$ws.Range("E46").Value = "  +1.74%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.68"
$ws.Range("E47").Value = "  -2.65%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.89"
$ws.Range("E48").Value = "  +11.07%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.58"
$ws.Range("E49").Value = "  -4.83%  "
$ws.Range("E50").Value = "  +0.16%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.98"
$ws.Range("E51").Value = "  +2.51%  "
